$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 268
    $ws.Range("F3").Value = 1384
    $ws.Range("F4").Value = 161
    $ws.Range("F6").Value = 235
    $ws.Range("F9").Value = 185
    $ws.Range("F11").Value = 4647
    $ws.Range("F12").Value = 6914
    $ws.Range("F16").Value = 573
    $ws.Range("F18").Value = 4149
    $ws.Range("F19").Value = 768
    $ws.Range("F22").Value = 2734
    $ws.Range("F24").Value = 550
    $ws.Range("F26").Value = 376
    $ws.Range("F27").Value = 374
    $ws.Range("F28").Value = 404
    $ws.Range("F29").Value = 233
    $ws.Range("F30").Value = 44
    $ws.Range("F31").Value = 1638
    $ws.Range("F32").Value = 1034
    $ws.Range("F34").Value = 397
    $ws.Range("F35").Value = 88
    $ws.Range("F36").Value = 550
    $ws.Range("F38").Value = 10
    $ws.Range("F40").Value = 171
    $ws.Range("F41").Value = 647
    $ws.Range("F42").Value = 19
}
